$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1  = 0.085453833382430844
    2  = -0.0099999996170581085
    3  = -0.0089999996211798106
    4  = 0.28399330306465131
    5  = -0.0059999996321282723
    6  = -0.0059999996184139093
    7  = -0.02041965385223854
    8  = -0.019999999543082403
    9  = -0.0059999996060602356
    10 = -0.0059999996010517975
    11 = -0.004499999607773475
    12 = -0.0059999995990240862
    13 = -0.0059999995939135076
    14 = -0.011999999562956276
    15 = -0.0059999995913777582
    16 = 0.001488976235449968
    17 = -0.0059999995887816127
    18 = -0.0089999995736445015
    19 = -0.0089999996214342737
    20 = -0.0089999996147831496
    21 = -0.034240814966644884
    22 = -0.0089999996130423199
    23 = -0.0089999996178482533
    24 = -0.041999999449615188
    25 = -0.041999999446694858
    26 = -0.0059999996172059866
    27 = -0.0059999996155859492
    28 = 0.011784857688669881
    29 = -0.011999999574221931
    30 = -0.019999999532218204
    31 = -0.016352708708746277
    32 = -0.020999999522738122
    33 = -0.0059999995962245478
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
